{"js": "// Updates the worksheet date header and all 100 addition/subtraction\n// problem cells in the table. Each entry is [oldText, newText] and the\n// entries are listed in document order (date paragraph first, then the\n// table cells row by row, left to right), matching the order returned by\n// `context.document.body.paragraphs`.\nconst replacements = [\n  [\"2024-03-07 Thursday\", \"2024-03-08 Friday\"],\n  [\"91-19=\", \"22+7=\"],\n  [\"38-1=\", \"96-57=\"],\n  [\"86-82=\", \"7+15=\"],\n  [\"45-25=\", \"65+8=\"],\n  [\"24+32=\", \"85-59=\"],\n  [\"54+3=\", \"76-27=\"],\n  [\"28+7=\", \"90-50=\"],\n  [\"25+70=\", \"85-9=\"],\n  [\"23+24=\", \"81-3=\"],\n  [\"40-16=\", \"39+59=\"],\n  [\"16+2=\", \"65-38=\"],\n  [\"83-15=\", \"87+3=\"],\n  [\"32+62=\", \"10+50=\"],\n  [\"26+21=\", \"99-63=\"],\n  [\"78-5=\", \"67-35=\"],\n  [\"29-27=\", \"22+28=\"],\n  [\"97-8=\", \"97-49=\"],\n  [\"48+22=\", \"19+58=\"],\n  [\"96-2=\", \"76-47=\"],\n  [\"46-8=\", \"50-5=\"],\n  [\"34-29=\", \"31-1=\"],\n  [\"59-28=\", \"29+35=\"],\n  [\"84-27=\", \"38-36=\"],\n  [\"8+52=\", \"64-17=\"],\n  [\"6+88=\", \"78+21=\"],\n  [\"22-3=\", \"12+31=\"],\n  [\"2+57=\", \"26+32=\"],\n  [\"8+30=\", \"69-14=\"],\n  [\"87-13=\", \"61-53=\"],\n  [\"59+20=\", \"90+0=\"],\n  [\"62-33=\", \"97-35=\"],\n  [\"28+27=\", \"5+42=\"],\n  [\"1+34=\", \"69-40=\"],\n  [\"80-78=\", \"13+46=\"],\n  [\"5+79=\", \"92-13=\"],\n  [\"29+30=\", \"40-9=\"],\n  [\"45+11=\", \"95-75=\"],\n  [\"50-39=\", \"93-26=\"],\n  [\"62+23=\", \"47+31=\"],\n  [\"66-33=\", \"36-5=\"],\n  [\"86-7=\", \"51+1=\"],\n  [\"74-1=\", \"75-20=\"],\n  [\"40-31=\", \"97-50=\"],\n  [\"91-47=\", \"85-57=\"],\n  [\"79-70=\", \"77+9=\"],\n  [\"27+22=\", \"42-28=\"],\n  [\"58+32=\", \"26-13=\"],\n  [\"28+46=\", \"49+32=\"],\n  [\"46+44=\", \"2+76=\"],\n  [\"73-26=\", \"15+81=\"],\n  [\"49+16=\", \"62-57=\"],\n  [\"27+58=\", \"87-87=\"],\n  [\"32+28=\", \"23-18=\"],\n  [\"36-32=\", \"82+4=\"],\n  [\"51-6=\", \"95-52=\"],\n  [\"63-3=\", \"52-23=\"],\n  [\"57+23=\", \"36+33=\"],\n  [\"92-79=\", \"97-58=\"],\n  [\"87-30=\", \"82-62=\"],\n  [\"56+5=\", \"44-15=\"],\n  [\"72-27=\", \"41-8=\"],\n  [\"68-62=\", \"43+32=\"],\n  [\"49+26=\", \"32+0=\"],\n  [\"34+59=\", \"90-88=\"],\n  [\"74-8=\", \"1+41=\"],\n  [\"85+13=\", \"27+34=\"],\n  [\"49-40=\", \"58-0=\"],\n  [\"7+4=\", \"90-14=\"],\n  [\"16+12=\", \"30-22=\"],\n  [\"23+61=\", \"91-58=\"],\n  [\"44-19=\", \"16+16=\"],\n  [\"28+2=\", \"24-16=\"],\n  [\"74-19=\", \"55-24=\"],\n  [\"95-50=\", \"31-17=\"],\n  [\"70-23=\", \"89-32=\"],\n  [\"7-4=\", \"45+44=\"],\n  [\"88-79=\", \"18+33=\"],\n  [\"10+33=\", \"65+14=\"],\n  [\"53-4=\", \"22+73=\"],\n  [\"34+60=\", \"90-64=\"],\n  [\"38+24=\", \"26+12=\"],\n  [\"85-78=\", \"45+33=\"],\n  [\"4+48=\", \"77-65=\"],\n  [\"56-10=\", \"8+87=\"],\n  [\"10+27=\", \"36-30=\"],\n  [\"72-48=\", \"57-49=\"],\n  [\"24+66=\", \"2+28=\"],\n  [\"45+15=\", \"75-8=\"],\n  [\"65-24=\", \"90-81=\"],\n  [\"78-33=\", \"49+3=\"],\n  [\"33+47=\", \"22-6=\"],\n  [\"67-11=\", \"35+14=\"],\n  [\"4+76=\", \"34-13=\"],\n  [\"40+47=\", \"0+38=\"],\n  [\"2+55=\", \"19-1=\"],\n  [\"59-42=\", \"46-26=\"],\n  [\"65-58=\", \"2+44=\"],\n  [\"65-46=\", \"99-83=\"],\n  [\"16+80=\", \"66+13=\"],\n  [\"66-55=\", \"10+46=\"]\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nif (items.length !== replacements.length) {\n  throw new Error(\n    \"Paragraph count mismatch: expected \" + replacements.length + \" but found \" + items.length\n  );\n}\n\nfor (let i = 0; i < items.length; i++) {\n  const oldText = replacements[i][0];\n  const newText = replacements[i][1];\n  const para = items[i];\n  if (para.text !== oldText) {\n    throw new Error(\n      \"Paragraph \" + i + \" text mismatch: expected \" + JSON.stringify(oldText) +\n      \" but found \" + JSON.stringify(para.text)\n    );\n  }\n  para.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Updates the worksheet date header and all 100 addition/subtraction problem\n# cells in the table. Each entry is (ParagraphIndex, oldText, newText). The\n# Word object model's $d.Paragraphs collection walks the document in order\n# (the date paragraph first, then the table cells row by row, left to\n# right) and additionally reports one extra empty paragraph at the end of\n# every table row (the row-end mark) that carries no text of its own -\n# those indices are simply skipped below.\n$replacements = @(\n    @(1, '2024-03-07 Thursday', '2024-03-08 Friday'),\n    @(2, '91-19=', '22+7='),\n    @(3, '38-1=', '96-57='),\n    @(4, '86-82=', '7+15='),\n    @(5, '45-25=', '65+8='),\n    @(6, '24+32=', '85-59='),\n    @(8, '54+3=', '76-27='),\n    @(9, '28+7=', '90-50='),\n    @(10, '25+70=', '85-9='),\n    @(11, '23+24=', '81-3='),\n    @(12, '40-16=', '39+59='),\n    @(14, '16+2=', '65-38='),\n    @(15, '83-15=', '87+3='),\n    @(16, '32+62=', '10+50='),\n    @(17, '26+21=', '99-63='),\n    @(18, '78-5=', '67-35='),\n    @(20, '29-27=', '22+28='),\n    @(21, '97-8=', '97-49='),\n    @(22, '48+22=', '19+58='),\n    @(23, '96-2=', '76-47='),\n    @(24, '46-8=', '50-5='),\n    @(26, '34-29=', '31-1='),\n    @(27, '59-28=', '29+35='),\n    @(28, '84-27=', '38-36='),\n    @(29, '8+52=', '64-17='),\n    @(30, '6+88=', '78+21='),\n    @(32, '22-3=', '12+31='),\n    @(33, '2+57=', '26+32='),\n    @(34, '8+30=', '69-14='),\n    @(35, '87-13=', '61-53='),\n    @(36, '59+20=', '90+0='),\n    @(38, '62-33=', '97-35='),\n    @(39, '28+27=', '5+42='),\n    @(40, '1+34=', '69-40='),\n    @(41, '80-78=', '13+46='),\n    @(42, '5+79=', '92-13='),\n    @(44, '29+30=', '40-9='),\n    @(45, '45+11=', '95-75='),\n    @(46, '50-39=', '93-26='),\n    @(47, '62+23=', '47+31='),\n    @(48, '66-33=', '36-5='),\n    @(50, '86-7=', '51+1='),\n    @(51, '74-1=', '75-20='),\n    @(52, '40-31=', '97-50='),\n    @(53, '91-47=', '85-57='),\n    @(54, '79-70=', '77+9='),\n    @(56, '27+22=', '42-28='),\n    @(57, '58+32=', '26-13='),\n    @(58, '28+46=', '49+32='),\n    @(59, '46+44=', '2+76='),\n    @(60, '73-26=', '15+81='),\n    @(62, '49+16=', '62-57='),\n    @(63, '27+58=', '87-87='),\n    @(64, '32+28=', '23-18='),\n    @(65, '36-32=', '82+4='),\n    @(66, '51-6=', '95-52='),\n    @(68, '63-3=', '52-23='),\n    @(69, '57+23=', '36+33='),\n    @(70, '92-79=', '97-58='),\n    @(71, '87-30=', '82-62='),\n    @(72, '56+5=', '44-15='),\n    @(74, '72-27=', '41-8='),\n    @(75, '68-62=', '43+32='),\n    @(76, '49+26=', '32+0='),\n    @(77, '34+59=', '90-88='),\n    @(78, '74-8=', '1+41='),\n    @(80, '85+13=', '27+34='),\n    @(81, '49-40=', '58-0='),\n    @(82, '7+4=', '90-14='),\n    @(83, '16+12=', '30-22='),\n    @(84, '23+61=', '91-58='),\n    @(86, '44-19=', '16+16='),\n    @(87, '28+2=', '24-16='),\n    @(88, '74-19=', '55-24='),\n    @(89, '95-50=', '31-17='),\n    @(90, '70-23=', '89-32='),\n    @(92, '7-4=', '45+44='),\n    @(93, '88-79=', '18+33='),\n    @(94, '10+33=', '65+14='),\n    @(95, '53-4=', '22+73='),\n    @(96, '34+60=', '90-64='),\n    @(98, '38+24=', '26+12='),\n    @(99, '85-78=', '45+33='),\n    @(100, '4+48=', '77-65='),\n    @(101, '56-10=', '8+87='),\n    @(102, '10+27=', '36-30='),\n    @(104, '72-48=', '57-49='),\n    @(105, '24+66=', '2+28='),\n    @(106, '45+15=', '75-8='),\n    @(107, '65-24=', '90-81='),\n    @(108, '78-33=', '49+3='),\n    @(110, '33+47=', '22-6='),\n    @(111, '67-11=', '35+14='),\n    @(112, '4+76=', '34-13='),\n    @(113, '40+47=', '0+38='),\n    @(114, '2+55=', '19-1='),\n    @(116, '59-42=', '46-26='),\n    @(117, '65-58=', '2+44='),\n    @(118, '65-46=', '99-83='),\n    @(119, '16+80=', '66+13='),\n    @(120, '66-55=', '10+46=')\n)\n\nforeach ($item in $replacements) {\n    $idx = $item[0]\n    $oldText = $item[1]\n    $newText = $item[2]\n\n    $p = $d.Paragraphs.Item($idx)\n    $r = $p.Range\n    # Range.Text includes the trailing paragraph/cell mark; trim it off so\n    # we compare against and overwrite only the visible text.\n    [void]$r.MoveEnd(1, -1)\n\n    if ($r.Text -ne $oldText) {\n        throw (\"Paragraph \" + $idx + \" text mismatch: expected '\" + $oldText + \"' but found '\" + $r.Text + \"'\")\n    }\n\n    $r.Text = $newText\n}\n"}
